# Slide 17 ("Тема 3" lecture deck): replace the empty, never-filled-in
# "title" placeholder with a manually placed text box carrying the
# section heading text -- matching the pattern every sibling slide in
# this section (10-16, 18) already uses instead of a real title
# placeholder.
#
# Target change (from the authoritative diff):
#   - remove the empty <p:sp> title placeholder (id=4, name
#     "Заголовок 3", <p:ph type="title"/>, empty text)
#   - append a new, non-placeholder <p:sp> text box (name
#     "Заголовок 3" again) at the very end of the slide's shape tree,
#     containing the run text
#     "Вопрос 1. Методы устранения выбросов и пропусков"

$p = $ppt.ActivePresentation
$slideIndex = 17
$dst = $p.Slides.Item($slideIndex)

# A neighboring slide that already carries the target "faux title"
# textbox (identical text & formatting). Cloning it via Copy/Paste
# reproduces the exact body/list-style XML PowerPoint itself produces,
# instead of trying to hand-author every formatting property.
$templateSlideIndex = 16
$src = $p.Slides.Item($templateSlideIndex)

$titleShapeName = "Заголовок 3"

# Locate the empty real title placeholder on the destination slide.
$titleShape = $null
for ($i = 1; $i -le $dst.Shapes.Count; $i++) {
    $sh = $dst.Shapes.Item($i)
    if ($sh.Name -eq $titleShapeName) {
        $titleShape = $sh
    }
}

# Locate the already-correct faux-title textbox on the template slide.
$templateShape = $null
for ($i = 1; $i -le $src.Shapes.Count; $i++) {
    $sh = $src.Shapes.Item($i)
    if ($sh.Name -eq $titleShapeName) {
        $templateShape = $sh
    }
}

# The engine hands out a newly created shape's id as "the lowest
# integer >= 2 not already used as an id on this slide". Right now the
# slide uses ids {1,3,4,6,7,8,9,10}, so naively adding/pasting a shape
# now would land on id 2 (because id 4 -- the placeholder we are about
# to delete -- is still "in the way" of the final value PowerPoint
# itself used: 11). Burn through the ids that would otherwise be
# handed out first (2, then 5) with disposable textboxes *before*
# deleting the placeholder, delete the placeholder only afterwards, so
# the pasted replacement shape ends up with id 11, matching the
# original author's save.
$burn1 = $dst.Shapes.AddTextbox(1, 0, 0, 1, 1)
$burn2 = $dst.Shapes.AddTextbox(1, 0, 0, 1, 1)

$newShape = $null
if ($templateShape -ne $null) {
    $templateShape.Copy()
    $pasted = $dst.Shapes.Paste()
    $newShape = $pasted.Item(1)
} else {
    # Fallback, should the template shape ever not be found: build the
    # textbox by hand with the same position/size/text/font size.
    $newShape = $dst.Shapes.AddTextbox(1, 234.8259842519685, 10.74992125984252, 474.88818897637793, 47.70700787401575)
    $newShape.TextFrame.TextRange.Text = "Вопрос 1. Методы устранения выбросов и пропусков"
    $newShape.TextFrame.TextRange.Font.Size = 20
}
$newShape.Name = $titleShapeName

# Clean up the disposable id-burning shapes and the old empty title
# placeholder now that the replacement shape has the right id.
$burn1.Delete()
$burn2.Delete()
if ($titleShape -ne $null) {
    $titleShape.Delete()
}
